$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Row 2
$ws.Range("D2").Value = 0.00008888868615031242
$ws.Range("E2").Value = 0.01573287695646286
$ws.Range("G2").Value = 0.001481673214584589
$ws.Range("H2").Value = 0.002649994101375341
$ws.Range("I2").Value = 0.004090444650501013
$ws.Range("J2").Value = 0.005886571016162634
$ws.Range("K2").Value = 0.0004930472932755947

# Row 3
$ws.Range("D3").Value = 0.001206391956657171
$ws.Range("E3").Value = 0.01717302855104208
$ws.Range("G3").Value = 0.001531261950731277
$ws.Range("H3").Value = 0.003925590310245752
$ws.Range("I3").Value = 0.003662615548819304
$ws.Range("J3").Value = 0.006333271041512489
$ws.Range("K3").Value = 0.0005031684413552284

# Row 4
$ws.Range("D4").Value = 0.0009107915684580803
$ws.Range("E4").Value = 0.01791709661483765
$ws.Range("G4").Value = 0.001692014280706644
$ws.Range("H4").Value = 0.003471673931926489
$ws.Range("I4").Value = 0.004269152414053679
$ws.Range("J4").Value = 0.006608132738620043
$ws.Range("K4").Value = 0.0005637500435113907

# Row 5
$ws.Range("D5").Value = 0.0004184572026133537
$ws.Range("E5").Value = 0.01623778697103262
$ws.Range("G5").Value = 0.001558565068989992
$ws.Range("H5").Value = 0.002985695842653513
$ws.Range("I5").Value = 0.004140559118241072
$ws.Range("J5").Value = 0.005899178795516491
$ws.Range("K5").Value = 0.0005072555504739285

# Row 6
$ws.Range("D6").Value = 0.001270269509404898
$ws.Range("E6").Value = 0.03900506906211376
$ws.Range("G6").Value = 0.002205895259976387
$ws.Range("H6").Value = 0.005577171687036753
$ws.Range("I6").Value = 0.02243523253127933
$ws.Range("J6").Value = 0.005772464442998171
$ws.Range("K6").Value = 0.0008078101091086864

# Row 8
$ws.Range("D8").Value = 0.00008888868615031242
$ws.Range("E8").Value = 0.01573287695646286
$ws.Range("G8").Value = 0.001481673214584589
$ws.Range("H8").Value = 0.002649994101375341
$ws.Range("I8").Value = 0.004090444650501013
$ws.Range("J8").Value = 0.005886571016162634
$ws.Range("K8").Value = 0.0004930472932755947

# Row 9
$ws.Range("D9").Value = 0.001206391956657171
$ws.Range("E9").Value = 0.01717302855104208
$ws.Range("G9").Value = 0.001531261950731277
$ws.Range("H9").Value = 0.003925590310245752
$ws.Range("I9").Value = 0.003662615548819304
$ws.Range("J9").Value = 0.006333271041512489
$ws.Range("K9").Value = 0.0005031684413552284

# Row 10
$ws.Range("D10").Value = 0.0009107915684580803
$ws.Range("E10").Value = 0.01791709661483765
$ws.Range("G10").Value = 0.001692014280706644
$ws.Range("H10").Value = 0.003471673931926489
$ws.Range("I10").Value = 0.004269152414053679
$ws.Range("J10").Value = 0.006608132738620043
$ws.Range("K10").Value = 0.0005637500435113907

# Row 11
$ws.Range("D11").Value = 0.0004184572026133537
$ws.Range("E11").Value = 0.01623778697103262
$ws.Range("G11").Value = 0.001558565068989992
$ws.Range("H11").Value = 0.002985695842653513
$ws.Range("I11").Value = 0.004140559118241072
$ws.Range("J11").Value = 0.005899178795516491
$ws.Range("K11").Value = 0.0005072555504739285

# Row 12
$ws.Range("D12").Value = 0.001270269509404898
$ws.Range("E12").Value = 0.03900506906211376
$ws.Range("G12").Value = 0.002205895259976387
$ws.Range("H12").Value = 0.005577171687036753
$ws.Range("I12").Value = 0.02243523253127933
$ws.Range("J12").Value = 0.005772464442998171
$ws.Range("K12").Value = 0.0008078101091086864
